# Insert a new data row at row 641 (pushes existing rows 641.. down by one)
# and populate it with a new weekly price observation, matching the
# existing "Vega Monumental Concepción - Limón" record layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 641, shifting rows 641:669 down to 642:670
$ws.Rows.Item(641).Insert()

# Populate the new row 641 with the new observation
$ws.Cells.Item(641, 1).Value = 11
$ws.Cells.Item(641, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(641, 3).Value = "Bíobío"
$ws.Cells.Item(641, 4).Value = (Get-Date -Year 2023 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(641, 5).Value = 8
$ws.Cells.Item(641, 6).Value = "Fruta"
$ws.Cells.Item(641, 7).Value = 100102
$ws.Cells.Item(641, 8).Value = "Cítricos"
$ws.Cells.Item(641, 9).Value = 100102003
$ws.Cells.Item(641, 10).Value = "Limón"
$ws.Cells.Item(641, 11).Value = "Sin especificar"
$ws.Cells.Item(641, 12).Value = "1a amarillo"
$ws.Cells.Item(641, 13).Value = 420
$ws.Cells.Item(641, 14).Value = 12000
$ws.Cells.Item(641, 15).Value = 14000
$ws.Cells.Item(641, 16).Value = 12952
$ws.Cells.Item(641, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(641, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(641, 19).Value = 810
$ws.Cells.Item(641, 20).Value = 16
